$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 14: turn the existing last row into a "border" (closing) row ---
# It currently uses the plain styles (B=4, C/D/E=5); the target reuses the
# bordered styles (B=6, C/D/E=7) and gains an empty styled A14 cell, same as
# rows 4 / 7 / 10 / 12 already do. Copy the formatting only (from row 12, an
# existing bordered row) so the style swap happens without touching the
# values already stored in B14:E14.
$ws.Range("A12:E12").Copy()
$ws.Range("A14:E14").PasteSpecial(-4122)

# --- Row 15: plain row (same family as rows 5/6, 8/9, 11) ---
$ws.Range("B11:E11").Copy()
$ws.Range("B15:E15").PasteSpecial(-4122)
$ws.Range("B15").Value = 240
$ws.Range("C15").Value = " Ah, welcome back! ♪"
$ws.Range("D15").Value = " Ах, с возвращением! ♪"
$ws.Range("E15").Value = " Àö, ò âïèâñàþåîéåí! ♪"

# --- Row 16: plain row, taller (matches rows 6/9/10's 21.6 family) ---
$ws.Range("B13:E13").Copy()
$ws.Range("B16:E16").PasteSpecial(-4122)
$ws.Rows.Item(16).RowHeight = 21.6
$ws.Range("B16").Value = 243
$ws.Range("C16").Value = " We were very worried about\nyou. No one seemed to know where you went!"
$ws.Range("D16").Value = " Мы за вас очень переживали.\nНикто не знал, куда вас занесло!"
$ws.Range("E16").Value = " Íú èà âàò ïœåîû ðåñåçéâàìé.\nÎéëóï îå èîàì, ëôäà âàò èàîåòìï!"

# --- Row 17: "border" row (closes the 15/16/17 conversation block) ---
$ws.Range("A12:E12").Copy()
$ws.Range("A17:E17").PasteSpecial(-4122)
$ws.Rows.Item(17).RowHeight = 31.8
$ws.Range("B17").Value = 246
$ws.Range("C17").Value = " It is good to see you home safe!\nAnd good to have you shop with us again! ♪"
$ws.Range("D17").Value = " Рад, что вы вернулись в целости\nи сохранности! И также, я рад, что вы\nснова можете вести с нами дела! ♪"
$ws.Range("E17").Value = " Ñàä, œóï âú âåñîôìéòû â øåìïòóé\né òïöñàîîïòóé! É óàëçå, ÿ ñàä, œóï âú\nòîïâà íïçåóå âåòóé ò îàíé äåìà! ♪"

# --- Row 18: plain row, taller (same family as row 16) ---
$ws.Range("B13:E13").Copy()
$ws.Range("B18:E18").PasteSpecial(-4122)
$ws.Rows.Item(18).RowHeight = 21.6
$ws.Range("B18").Value = 202
$ws.Range("C18").Value = " My thoughts go with you! I do\nhope the best for you!"
$ws.Range("D18").Value = " Все мои мысли только о вас!\nЯ желаю вам всего наилучшего!"
$ws.Range("E18").Value = " Âòå íïé íúòìé óïìûëï ï âàò!\nŸ çåìàý âàí âòåãï îàéìôœšåãï!"

# --- sheet view: scroll / selection mirror the author's final state ---
$ws.Range("A13").Select()
$ws.Range("D18").Select()

$excel.ActiveWindow.ScrollRow = 13
